# Added games for 1/17/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing "Beat Vegas?" (column G) values for the existing
#     1/16/2021 games (rows 79-84) that were left blank. ---
$ws.Range("G79").Value = "No"
$ws.Range("G80").Value = "No"
$ws.Range("G81").Value = "Yes"
$ws.Range("G82").Value = "No"
$ws.Range("G83").Value = "No"
$ws.Range("G84").Value = "Yes"

# --- Append the new games played on 1/17/2021 (serial date 44213) ---
$newGames = @(
    @{ Row = 85; Home = "BOS"; Away = "NYK"; Spread = -7;   Pred = -18;   Diff = 11 },
    @{ Row = 86; Home = "DAL"; Away = "CHI"; Spread = -7;   Pred = 13.8;  Diff = -20.8 },
    @{ Row = 87; Home = "OKC"; Away = "PHI"; Spread = 2.5;  Pred = 3.7;   Diff = -1.2 },
    @{ Row = 88; Home = "DEN"; Away = "UTA"; Spread = 1;    Pred = 0.3;   Diff = 0.7 },
    @{ Row = 89; Home = "SAC"; Away = "NOP"; Spread = 2.5;  Pred = -10.7; Diff = 13.2 },
    @{ Row = 90; Home = "LAC"; Away = "IND"; Spread = -6.5; Pred = -6;    Diff = -0.5 }
)

foreach ($game in $newGames) {
    $r = $game.Row
    $ws.Range("A$r").Value = 44213
    $ws.Range("A$r").NumberFormat = "yyyy\-mm\-dd"
    $ws.Range("B$r").Value = $game.Home
    $ws.Range("C$r").Value = $game.Away
    $ws.Range("D$r").Value = $game.Spread
    $ws.Range("E$r").Value = $game.Pred
    $ws.Range("F$r").Value = $game.Diff
}

# Match the saved selection state recorded in the workbook after the edit.
$ws.Range("H86").Select()
